$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 0.671980650564856
$ws.Range("C8").Value = 28.9076813004957
$ws.Range("D8").Value = 54.92132823176864
$ws.Range("E8").Value = 4331.448889646987
$ws.Range("F8").Value = 479389.9575332087
$ws.Range("G8").Value = 24942430.48433685
$ws.Range("H8").Value = 1523777560.795277
$ws.Range("I8").Value = 53371990299098.23
